$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 222; this shifts the existing row 222..294
# block down to 223..295 and keeps all their values/formatting intact.
$ws.Rows.Item(222).Insert()

# Populate the newly inserted row 222 with a new price observation.
# Columns that stay constant across this block (A,B,C,E,F,G,H,I,N,Q,R) are
# copied from the row immediately below (the row that used to be 222).
$valA = $ws.Cells.Item(223, 1).Value2
$valB = $ws.Cells.Item(223, 2).Value2
$valC = $ws.Cells.Item(223, 3).Value2
$valE = $ws.Cells.Item(223, 5).Value2
$valF = $ws.Cells.Item(223, 6).Value2
$valG = $ws.Cells.Item(223, 7).Value2
$valH = $ws.Cells.Item(223, 8).Value2
$valI = $ws.Cells.Item(223, 9).Value2
$valN = $ws.Cells.Item(223, 14).Value2
$valQ = $ws.Cells.Item(223, 17).Value2
$valR = $ws.Cells.Item(223, 18).Value2

$ws.Cells.Item(222, 1).Value = $valA
$ws.Cells.Item(222, 2).Value = $valB
$ws.Cells.Item(222, 3).Value = $valC
$ws.Cells.Item(222, 4).Value = 44876
$ws.Cells.Item(222, 5).Value = $valE
$ws.Cells.Item(222, 6).Value = $valF
$ws.Cells.Item(222, 7).Value = $valG
$ws.Cells.Item(222, 8).Value = $valH
$ws.Cells.Item(222, 9).Value = $valI
$ws.Cells.Item(222, 10).Value = 75
$ws.Cells.Item(222, 11).Value = 8000
$ws.Cells.Item(222, 12).Value = 9000
$ws.Cells.Item(222, 13).Value = 8467
$ws.Cells.Item(222, 14).Value = $valN
$ws.Cells.Item(222, 15).Value = "Región Metropolitana"
$ws.Cells.Item(222, 16).Value = 339
$ws.Cells.Item(222, 17).Value = $valQ
$ws.Cells.Item(222, 18).Value = $valR
